# Apply the "Hong Kong EPS v2.0.0" update to the SoCDTtiNTY workbook.

$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("SoCDTtiNTY")

# --- Update the data values on the SoCDTtiNTY sheet ---------------------
$wsData.Range("B2").Value = 0.076
$wsData.Range("C2").Value = 0.07

$wsData.Range("B3").Value = 0.0435
$wsData.Range("C3").Value = 0.035

$wsData.Range("C4").Value = 0.042

$wsData.Range("B5").Value = 0.029

$wsData.Range("B7").Value = 0.0587

# --- Add the new header cell / label in A1 -------------------------------
$wsData.Range("A1").Value = "Share that is New (dimensionless)"
$wsData.Range("A1").Font.Bold = $true
$wsData.Range("A1").WrapText = $true
$wsData.Rows.Item(1).RowHeight = 60

# --- Switch which sheet/tab is active ------------------------------------
$wsAbout.Activate()
